$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 301-302; this shifts the existing rows 301:330
# down to 303:332, exactly like pressing "Insert Rows" in Excel.
$ws.Rows("301:302").Insert()

# Row 301 - new weekly price entry (Primera)
$ws.Range("A301").Value = 11
$ws.Range("B301").Value = "Vega Monumental Concepción"
$ws.Range("C301").Value = "Bíobío"
$ws.Range("D301").Value = 45142
$ws.Range("E301").Value = 8
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100101
$ws.Range("H301").Value = "Berries"
$ws.Range("I301").Value = 100101007
$ws.Range("J301").Value = "Kiwi"
$ws.Range("K301").Value = "Hayward"
$ws.Range("L301").Value = "Primera"
$ws.Range("M301").Value = 100
$ws.Range("N301").Value = 14000
$ws.Range("O301").Value = 15000
$ws.Range("P301").Value = 14500
$ws.Range("Q301").Value = "$/bandeja 18 kilos"
$ws.Range("R301").Value = "Región de O'Higgins"
$ws.Range("S301").Value = 806
$ws.Range("T301").Value = 18

# Row 302 - new weekly price entry (Segunda)
$ws.Range("A302").Value = 11
$ws.Range("B302").Value = "Vega Monumental Concepción"
$ws.Range("C302").Value = "Bíobío"
$ws.Range("D302").Value = 45142
$ws.Range("E302").Value = 8
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100101
$ws.Range("H302").Value = "Berries"
$ws.Range("I302").Value = 100101007
$ws.Range("J302").Value = "Kiwi"
$ws.Range("K302").Value = "Hayward"
$ws.Range("L302").Value = "Segunda"
$ws.Range("M302").Value = 50
$ws.Range("N302").Value = 12000
$ws.Range("O302").Value = 12000
$ws.Range("P302").Value = 12000
$ws.Range("Q302").Value = "$/bandeja 18 kilos"
$ws.Range("R302").Value = "Región de O'Higgins"
$ws.Range("S302").Value = 667
$ws.Range("T302").Value = 18
